$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 52: "Press Pause" / "Pause/Running should be paused."
$ws.Range("B52").Value = "Press Pause"
$ws.Range("C52").Value = "Pause/Running should be paused."

# New row 53: "Press Run" / "Test execution should continue from the current target."
$ws.Range("B53").Value = "Press Run"
$ws.Range("C53").Value = "Test execution should continue from the current target."

# Apply same style (wrap text, vertical top) as the rest of column B/C body rows
$ws.Range("B52:C53").WrapText = $true
$ws.Range("B52:C53").VerticalAlignment = -4160

# Update selection / view to match new last cell
$ws.Range("C53").Select()
